# Track shipment history API implemented
# Applies the changes described by the commit:
#  - CRUD Branch (row 2): Status column now notes SP created / Java integration pending
#  - CRUD Customer (row 5): marked Complete, highlight fill removed (no longer "in progress")
#  - Track Shipment (row 13): marked Complete
#  - Column F (Status) widened to fit new text, row 2 height grown for wrapped text
#  - Selection left on F2 (Status cell for the newly-updated CRUD Branch row), scrolled to top

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- CRUD Branch (row 2): add status note ---
$ws.Range("F2").Value = "SP Created; `nIntegration with Java code pending"
$ws.Range("F2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 43.2

# --- CRUD Customer (row 5): complete + drop the "in-progress" highlight fill ---
$ws.Range("F5").Value = "Complete"
$ws.Range("B5:G5").Interior.Pattern = -4142   # xlPatternNone

# --- Track Shipment (row 13): complete ---
$ws.Range("F13").Value = "Complete"

# --- Column F (Status) widen to fit the longer note ---
$ws.Columns.Item(6).ColumnWidth = 15.88671875

# --- Final selection / scroll position ---
$ws.Range("F2").Select()
